$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update inserts a new data row for this market right after the
# existing row 83 (row 84 shifts down to 85, ..., old row 132 becomes 133).
# The newly inserted row 84 duplicates the content of row 83, matching the
# same weekly-refresh pattern used across this price series.
$ws.Rows.Item(83).Copy()
$ws.Rows.Item(84).Insert()
